{"js": "// Correct the exercise reference in the Overview paragraph:\n// \"(from the previous lab and from exercise 1.5)\"\n//   -> \"(from the previous lab and from exercise 2.B or 2.11)\"\nconst body = context.document.body;\nconst results = body.search(\"1.5\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"2.B or 2.11\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Correct the exercise reference in the Overview paragraph:\n# \"(from the previous lab and from exercise 1.5)\"\n#   -> \"(from the previous lab and from exercise 2.B or 2.11)\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"1.5\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"2.B or 2.11\"\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute(\n    \"1.5\",          # FindText\n    $false,         # MatchCase\n    $false,         # MatchWholeWord\n    $false,         # MatchWildcards\n    $false,         # MatchSoundsLike\n    $false,         # MatchAllWordForms\n    $true,          # Forward\n    1,              # Wrap (wdFindContinue)\n    $false,         # Format\n    \"2.B or 2.11\",  # ReplaceWith\n    2               # Replace (wdReplaceAll)\n)\n"}
